# Updates the per-line loading-percent results (Sheet1) to the recomputed
# values for the 380 kV case (Case_5_153/res_line/loading_percent.xlsx).
# Columns B-F and K-N of data rows 2-25 change; everything else is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 index 0)
$ws.Range("B2").Value = 18.5679874210492
$ws.Range("C2").Value = 5.562401437147237
$ws.Range("D2").Value = 8.904495809317439
$ws.Range("E2").Value = 10.67598311675268
$ws.Range("F2").Value = 43.52110274804083
$ws.Range("K2").Value = 14.91549994552527
$ws.Range("L2").Value = 10.11849393003555
$ws.Range("M2").Value = 16.99277315757985
$ws.Range("N2").Value = 24.89438636514322
# Row 3 (A3 index 1)
$ws.Range("B3").Value = 18.42006533862371
$ws.Range("C3").Value = 5.399777418781163
$ws.Range("D3").Value = 8.913683858099608
$ws.Range("E3").Value = 10.69221696583118
$ws.Range("F3").Value = 43.39044630862723
$ws.Range("K3").Value = 14.81024531050834
$ws.Range("L3").Value = 10.12974653347598
$ws.Range("M3").Value = 16.98428539354114
$ws.Range("N3").Value = 24.92284110228973
# Row 4 (A4 index 2)
$ws.Range("B4").Value = 18.33363894517318
$ws.Range("C4").Value = 5.296116411645936
$ws.Range("D4").Value = 8.919464365204497
$ws.Range("E4").Value = 10.70309845947502
$ws.Range("F4").Value = 43.31839245167132
$ws.Range("K4").Value = 14.74919779977436
$ws.Range("L4").Value = 10.13806042733651
$ws.Range("M4").Value = 16.98224569073491
$ws.Range("N4").Value = 24.94212439181819
# Row 5 (A5 index 3)
$ws.Range("B5").Value = 18.29956018811473
$ws.Range("C5").Value = 5.25294532627352
$ws.Range("D5").Value = 8.921855124916892
$ws.Range("E5").Value = 10.70776292067703
$ws.Range("F5").Value = 43.29109947485035
$ws.Range("K5").Value = 14.72524368886243
$ws.Range("L5").Value = 10.14180184388389
$ws.Range("M5").Value = 16.98221387853335
$ws.Range("N5").Value = 24.95043806957895
# Row 6 (A6 index 4)
$ws.Range("B6").Value = 18.2939713023221
$ws.Range("C6").Value = 5.245721643515959
$ws.Range("D6").Value = 8.922254237969398
$ws.Range("E6").Value = 10.70855136309504
$ws.Range("F6").Value = 43.28669289571229
$ws.Range("K6").Value = 14.72132253397177
$ws.Range("L6").Value = 10.14244445485612
$ws.Range("M6").Value = 16.9822569214856
$ws.Range("N6").Value = 24.95184606373742
# Row 7 (A7 index 5)
$ws.Range("B7").Value = 18.33317468432132
$ws.Range("C7").Value = 5.295537909298345
$ws.Range("D7").Value = 8.919496465231243
$ws.Range("E7").Value = 10.70316043360824
$ws.Range("F7").Value = 43.31801597110368
$ws.Range("K7").Value = 14.74887097868076
$ws.Range("L7").Value = 10.13810945410836
$ws.Range("M7").Value = 16.98224202297641
$ws.Range("N7").Value = 24.94223466848532
# Row 8 (A8 index 6)
$ws.Range("B8").Value = 18.51609248382123
$ws.Range("C8").Value = 5.507136106533902
$ws.Range("D8").Value = 8.907635120370385
$ws.Range("E8").Value = 10.68139110247534
$ws.Range("F8").Value = 43.47436652508966
$ws.Range("K8").Value = 14.87847970087586
$ws.Range("L8").Value = 10.12208240659655
$ws.Range("M8").Value = 16.98918946595218
$ws.Range("N8").Value = 24.9038214539037
# Row 9 (A9 index 7)
$ws.Range("B9").Value = 18.9080673296091
$ws.Range("C9").Value = 5.890700165091569
$ws.Range("D9").Value = 8.885468287321101
$ws.Range("E9").Value = 10.64593616881396
$ws.Range("F9").Value = 43.8450357607373
$ws.Range("K9").Value = 15.1599449718986
$ws.Range("L9").Value = 10.10178950066164
$ws.Range("M9").Value = 17.02788550263548
$ws.Range("N9").Value = 24.84287339385074
# Row 10 (A10 index 8)
$ws.Range("B10").Value = 19.21403278915717
$ws.Range("C10").Value = 6.152014953496629
$ws.Range("D10").Value = 8.869834362736025
$ws.Range("E10").Value = 10.6242759495916
$ws.Range("F10").Value = 44.15528379456642
$ws.Range("K10").Value = 15.381814092631
$ws.Range("L10").Value = 10.09365403966956
$ws.Range("M10").Value = 17.0714500490097
$ws.Range("N10").Value = 24.80686640754913
# Row 11 (A11 index 9)
$ws.Range("B11").Value = 19.35660910707412
$ws.Range("C11").Value = 6.26620596072532
$ws.Range("D11").Value = 8.862860524544434
$ws.Range("E11").Value = 10.61537040809885
$ws.Range("F11").Value = 44.30437325575956
$ws.Range("K11").Value = 15.48566327898215
$ws.Range("L11").Value = 10.09141953791619
$ws.Range("M11").Value = 17.0945142897878
$ws.Range("N11").Value = 24.79239165674452
# Row 12 (A12 index 10)
$ws.Range("B12").Value = 19.411041585999
$ws.Range("C12").Value = 6.308756699413584
$ws.Range("D12").Value = 8.860239344278318
$ws.Range("E12").Value = 10.61213402366852
$ws.Range("F12").Value = 44.36194650220674
$ws.Range("K12").Value = 15.525376659747
$ws.Range("L12").Value = 10.09078376085402
$ws.Range("M12").Value = 17.10371074812145
$ws.Range("N12").Value = 24.78718446527014
# Row 13 (A13 index 11)
$ws.Range("B13").Value = 19.39929965718643
$ws.Range("C13").Value = 6.299623623679058
$ws.Range("D13").Value = 8.860802991098318
$ws.Range("E13").Value = 10.61282499598253
$ws.Range("F13").Value = 44.34949788648812
$ws.Range("K13").Value = 15.51680693345438
$ws.Range("L13").Value = 10.0909113384434
$ws.Range("M13").Value = 17.10170962835011
$ws.Range("N13").Value = 24.78829373498687
# Row 14 (A14 index 12)
$ws.Range("B14").Value = 19.3610787070438
$ws.Range("C14").Value = 6.269720576118419
$ws.Range("D14").Value = 8.862644485577119
$ws.Range("E14").Value = 10.61510142611235
$ws.Range("F14").Value = 44.30908763436934
$ws.Range("K14").Value = 15.48892293275062
$ws.Range("L14").Value = 10.09136301880363
$ws.Range("M14").Value = 17.0952616443814
$ws.Range("N14").Value = 24.79195776414601
# Row 15 (A15 index 13)
$ws.Range("B15").Value = 19.33772345572406
$ws.Range("C15").Value = 6.251313631087253
$ws.Range("D15").Value = 8.863775008619998
$ws.Range("E15").Value = 10.61651349992694
$ws.Range("F15").Value = 44.28447970929592
$ws.Range("K15").Value = 15.47189274115141
$ws.Range("L15").Value = 10.09166706841934
$ws.Range("M15").Value = 17.09137216781336
$ws.Range("N15").Value = 24.79423778560099
# Row 16 (A16 index 14)
$ws.Range("B16").Value = 19.20477943593281
$ws.Range("C16").Value = 6.144456529471705
$ws.Range("D16").Value = 8.870292881589773
$ws.Range("E16").Value = 10.62487699069295
$ws.Range("F16").Value = 44.14569821113651
$ws.Range("K16").Value = 15.37508334511008
$ws.Range("L16").Value = 10.09382953692748
$ws.Range("M16").Value = 17.07000771589228
$ws.Range("N16").Value = 24.80785071271173
# Row 17 (A17 index 15)
$ws.Range("B17").Value = 19.12405881652843
$ws.Range("C17").Value = 6.077690702880411
$ws.Range("D17").Value = 8.87432661757156
$ws.Range("E17").Value = 10.63025023355347
$ws.Range("F17").Value = 44.06257881185376
$ws.Range("K17").Value = 15.31641926179657
$ws.Range("L17").Value = 10.09553138990165
$ws.Range("M17").Value = 17.05772989457807
$ws.Range("N17").Value = 24.81668982590669
# Row 18 (A18 index 16)
$ws.Range("B18").Value = 19.07795242391687
$ws.Range("C18").Value = 6.038848950275991
$ws.Range("D18").Value = 8.876659727571209
$ws.Range("E18").Value = 10.63343001284739
$ws.Range("F18").Value = 44.01552137454561
$ws.Range("K18").Value = 15.28295401885317
$ws.Range("L18").Value = 10.09664830434403
$ws.Range("M18").Value = 17.05097387928081
$ws.Range("N18").Value = 24.82195315023645
# Row 19 (A19 index 17)
$ws.Range("B19").Value = 19.06239827272146
$ws.Range("C19").Value = 6.025622831770979
$ws.Range("D19").Value = 8.877451919377819
$ws.Range("E19").Value = 10.63452196754765
$ws.Range("F19").Value = 43.99971827507944
$ws.Range("K19").Value = 15.27167176213633
$ws.Range("L19").Value = 10.09705019679844
$ws.Range("M19").Value = 17.04873906510517
$ws.Range("N19").Value = 24.82376601184039
# Row 20 (A20 index 18)
$ws.Range("B20").Value = 19.13261867052912
$ws.Range("C20").Value = 6.084843692202145
$ws.Range("D20").Value = 8.873895874295615
$ws.Range("E20").Value = 10.62966900996438
$ws.Range("F20").Value = 44.07134952088349
$ws.Range("K20").Value = 15.32263573644914
$ws.Range("L20").Value = 10.0953359400163
$ws.Range("M20").Value = 17.05900526390973
$ws.Range("N20").Value = 24.81573032936392
# Row 21 (A21 index 19)
$ws.Range("B21").Value = 19.37229350019709
$ws.Range("C21").Value = 6.278522707416094
$ws.Range("D21").Value = 8.862103062008998
$ws.Range("E21").Value = 10.61442909647135
$ws.Range("F21").Value = 44.32092703889811
$ws.Range("K21").Value = 15.49710286840493
$ws.Range("L21").Value = 10.09122464382396
$ws.Range("M21").Value = 17.09714305713095
$ws.Range("N21").Value = 24.79087411032312
# Row 22 (A22 index 20)
$ws.Range("B22").Value = 19.53148915514736
$ws.Range("C22").Value = 6.401069196810738
$ws.Range("D22").Value = 8.854510283510706
$ws.Range("E22").Value = 10.60526119724485
$ws.Range("F22").Value = 44.49053293958847
$ws.Range("K22").Value = 15.61337326886036
$ws.Range("L22").Value = 10.08976366244026
$ws.Range("M22").Value = 17.12476199022956
$ws.Range("N22").Value = 24.77622675897768
# Row 23 (A23 index 21)
$ws.Range("B23").Value = 19.44630496032127
$ws.Range("C23").Value = 6.336038225318778
$ws.Range("D23").Value = 8.85855228098881
$ws.Range("E23").Value = 10.6100818988508
$ws.Range("F23").Value = 44.39942663380306
$ws.Range("K23").Value = 15.55112275419752
$ws.Range("L23").Value = 10.09043141321245
$ws.Range("M23").Value = 17.10977630994208
$ws.Range("N23").Value = 24.78389809742712
# Row 24 (A24 index 22)
$ws.Range("B24").Value = 19.12874781959426
$ws.Range("C24").Value = 6.08161124845619
$ws.Range("D24").Value = 8.874090569502856
$ws.Range("E24").Value = 10.62993149882839
$ws.Range("F24").Value = 44.06738201040613
$ws.Range("K24").Value = 15.31982445123879
$ws.Range("L24").Value = 10.09542387147965
$ws.Range("M24").Value = 17.05842772634185
$ws.Range("N24").Value = 24.81616355206219
# Row 25 (A25 index 23)
$ws.Range("B25").Value = 18.79869717418557
$ws.Range("C25").Value = 5.790445355787345
$ws.Range("D25").Value = 8.891349454603443
$ws.Range("E25").Value = 10.6547554052835
$ws.Range("F25").Value = 43.73801068769019
$ws.Range("K25").Value = 15.08103502444421
$ws.Range("L25").Value = 10.10608839980877
$ws.Range("M25").Value = 17.01474682229501
$ws.Range("N25").Value = 24.85782142415339
